$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4146.763
$ws.Range("I15").Value = 4146.763
$ws.Range("K15").Value = 12440.289
$ws.Range("M15").Value = -12271.289
$ws.Range("H32").Value = 1427.7
$ws.Range("I32").Value = 3500
$ws.Range("J32").Value = 909.625
$ws.Range("K32").Value = 3500
$ws.Range("L32").Value = 909.625
$ws.Range("M32").Value = -3174
$ws.Range("N32").Value = -1561.625
$ws.Range("H33").Value = 117.8
$ws.Range("J33").Value = 122.71429
$ws.Range("L33").Value = 122.71429
$ws.Range("N33").Value = -580.71429
$ws.Range("H43").Value = 1200.7858
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1200.7858
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1200.7858
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1338.7858
$ws.Range("H62").Value = 1200
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1400
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1400
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -2648
$ws.Range("H65").Value = 1200
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1400
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -13240
$ws.Range("H98").Value = 4234.2144
$ws.Range("I98").Value = 3725.4546
$ws.Range("K98").Value = 3725.4546
$ws.Range("M98").Value = -2227.4546
$ws.Range("H100").Value = 1279.125
$ws.Range("J100").Value = 2195
$ws.Range("L100").Value = 2195
$ws.Range("N100").Value = -3277
$ws.Range("H111").Value = 1412.6
$ws.Range("I111").Value = 733.3333
$ws.Range("J111").Value = 2431.5
$ws.Range("K111").Value = 2199.9999
$ws.Range("L111").Value = 7294.5
$ws.Range("M111").Value = 867.0001000000002
$ws.Range("N111").Value = -13428.5
$ws.Range("H113").Value = 141201
$ws.Range("I113").Value = 141201
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 141201
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -137947
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 4234.2144
$ws.Range("I122").Value = 3725.4546
$ws.Range("K122").Value = 11176.3638
$ws.Range("M122").Value = -8726.363799999999
$ws.Range("H138").Value = 1570.375
$ws.Range("I138").Value = 1249.5161
$ws.Range("J138").Value = 2675.5557
$ws.Range("K138").Value = 3748.5483
$ws.Range("L138").Value = 8026.6671
$ws.Range("M138").Value = 1391.4517
$ws.Range("N138").Value = -18306.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1676.2667
$ws.Range("I61").Value = 1004.73914
$ws.Range("J61").Value = 3882.7144
$ws.Range("K61").Value = 1004.73914
$ws.Range("L61").Value = 3882.7144
$ws.Range("M61").Value = -792.73914
$ws.Range("N61").Value = -4306.7144
$ws.Range("H74").Value = 1383.8096
$ws.Range("I74").Value = 1149
$ws.Range("K74").Value = 1149
$ws.Range("M74").Value = -275
$ws.Range("H77").Value = 1383.8096
$ws.Range("I77").Value = 1149
$ws.Range("K77").Value = 5745
$ws.Range("M77").Value = -1377
$ws.Range("H110").Value = 1748.9333
$ws.Range("I110").Value = 1517.5385
$ws.Range("K110").Value = 1517.5385
$ws.Range("M110").Value = 527.4614999999999
$ws.Range("H132").Value = 1379.9767
$ws.Range("I132").Value = 1148.4242
$ws.Range("K132").Value = 3445.2726
$ws.Range("M132").Value = -915.2725999999998
$ws.Range("H136").Value = 1676.2667
$ws.Range("I136").Value = 1004.73914
$ws.Range("J136").Value = 3882.7144
$ws.Range("K136").Value = 3014.21742
$ws.Range("L136").Value = 11648.1432
$ws.Range("M136").Value = -464.2174199999999
$ws.Range("N136").Value = -16748.1432
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1117.5883
$ws.Range("I22").Value = 466.55554
$ws.Range("J22").Value = 1850
$ws.Range("K22").Value = 466.55554
$ws.Range("L22").Value = 1850
$ws.Range("M22").Value = -116.55554
$ws.Range("N22").Value = -2550
$ws.Range("H86").Value = 200002050
$ws.Range("I86").Value = 250001470
$ws.Range("K86").Value = 250001470
$ws.Range("M86").Value = -250000347
$ws.Range("H89").Value = 200002050
$ws.Range("I89").Value = 250001470
$ws.Range("K89").Value = 1250007350
$ws.Range("M89").Value = -1250001734
$ws.Range("H99").Value = 1999.5
$ws.Range("I99").Value = 1999.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1999.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -501.5
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 1999.5
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1386.659
$ws.Range("I132").Value = 942.46875
$ws.Range("K132").Value = 2827.40625
$ws.Range("M132").Value = -297.40625
$ws.Range("H134").Value = 1491.2322
$ws.Range("I134").Value = 1345.2727
$ws.Range("K134").Value = 4035.8181
$ws.Range("M134").Value = -1500.8181
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 15004.167
$ws.Range("I87").Value = 10006.25
$ws.Range("K87").Value = 30018.75
$ws.Range("M87").Value = -28770.75
$ws.Range("H90").Value = 15004.167
$ws.Range("I90").Value = 10006.25
$ws.Range("K90").Value = 90056.25
$ws.Range("M90").Value = -83816.25
$ws.Range("H92").Value = 318.375
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 292.42856
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 877.28568
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -3373.28568
$ws.Range("H114").Value = 2031.8
$ws.Range("I114").Value = 482
$ws.Range("J114").Value = 3065
$ws.Range("K114").Value = 1446
$ws.Range("L114").Value = 9195
$ws.Range("M114").Value = 1808
$ws.Range("N114").Value = -15703
$ws.Range("H131").Value = 7587208
$ws.Range("J131").Value = 14411.115
$ws.Range("L131").Value = 43233.345
$ws.Range("N131").Value = -53313.345
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3856.5715
$ws.Range("J80").Value = 4166
$ws.Range("L80").Value = 4166
$ws.Range("N80").Value = -6162
$ws.Range("H83").Value = 3856.5715
$ws.Range("J83").Value = 4166
$ws.Range("L83").Value = 20830
$ws.Range("N83").Value = -30814
$ws.Range("H113").Value = 1305
$ws.Range("I113").Value = 1332.5
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 1332.5
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = 837.5
$ws.Range("N113").Value = -5590
$ws.Range("H132").Value = 1834228.8
$ws.Range("I132").Value = 2566175.5
$ws.Range("K132").Value = 7698526.5
$ws.Range("M132").Value = -7695996.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3054.2856
$ws.Range("J22").Value = 1896.6666
$ws.Range("L22").Value = 1896.6666
$ws.Range("N22").Value = -2486.6666
$ws.Range("H27").Value = 3054.2856
$ws.Range("J27").Value = 1896.6666
$ws.Range("L27").Value = 1896.6666
$ws.Range("N27").Value = -2110.6666
$ws.Range("H46").Value = 1263.4
$ws.Range("I46").Value = 654.2
$ws.Range("K46").Value = 654.2
$ws.Range("M46").Value = -466.2
$ws.Range("H61").Value = 2346.2942
$ws.Range("I61").Value = 1992.4667
$ws.Range("K61").Value = 1992.4667
$ws.Range("M61").Value = -1790.4667
$ws.Range("H93").Value = 1096.5
$ws.Range("J93").Value = 1737.2
$ws.Range("L93").Value = 1737.2
$ws.Range("N93").Value = -4233.2
$ws.Range("H113").Value = 2346.2942
$ws.Range("I113").Value = 1992.4667
$ws.Range("K113").Value = 1992.4667
$ws.Range("M113").Value = 177.5333000000001
$ws.Range("H132").Value = 1221.6031
$ws.Range("I132").Value = 1007.6445
$ws.Range("J132").Value = 1756.5
$ws.Range("K132").Value = 3022.9335
$ws.Range("L132").Value = 5269.5
$ws.Range("M132").Value = -492.9335000000001
$ws.Range("N132").Value = -10329.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5474.4614
$ws.Range("J14").Value = 5113.7617
$ws.Range("L14").Value = 5113.7617
$ws.Range("N14").Value = -5449.7617
$ws.Range("H58").Value = 22531
$ws.Range("I58").Value = 5085
$ws.Range("J58").Value = 39977
$ws.Range("K58").Value = 5085
$ws.Range("L58").Value = 39977
$ws.Range("M58").Value = -4777
$ws.Range("N58").Value = -40593
$ws.Range("H96").Value = 9289.700000000001
$ws.Range("I96").Value = 2259.4
$ws.Range("K96").Value = 2259.4
$ws.Range("M96").Value = -886.4000000000001
$ws.Range("H107").Value = 557.8095
$ws.Range("I107").Value = 231.38461
$ws.Range("K107").Value = 694.15383
$ws.Range("M107").Value = 1225.84617
$ws.Range("H113").Value = 615.1875
$ws.Range("I113").Value = 295.69232
$ws.Range("K113").Value = 887.07696
$ws.Range("M113").Value = 1282.92304
$ws.Range("H136").Value = 16341699
$ws.Range("I136").Value = 18520126
$ws.Range("K136").Value = 55560378
$ws.Range("M136").Value = -55557828
